$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 <- source row 57
$ws.Range("B42").Value = "캥거 장스탠드 만달라키 노을 석양 조명 선셋 무드등 MEL001"
$ws.Range("C42").Value = "https://search.shopping.naver.com/gate.nhn?id=28236356554"
$ws.Range("D42").Value = "https://shopping-phinf.pstatic.net/main_2823635/28236356554.20210801004606.jpg"
$ws.Range("E42").Value = "22000"
$ws.Range("G42").Value = "네이버"
$ws.Range("H42").Value = "일반 - 가격비교 상품"
$ws.Range("I42").Value = "캥거"
$ws.Range("J42").Value = "델로나"
$ws.Range("N42").Value = "장스탠드"

# Row 43 <- source row 42
$ws.Range("B43").Value = "LED 코너 플로어 램프 거실 서 장식 분위기 램프 침실 장식 스탠드 조명 실내 조명"
$ws.Range("C43").Value = "https://search.shopping.naver.com/gate.nhn?id=24687593897"
$ws.Range("D43").Value = "https://shopping-phinf.pstatic.net/main_2468759/24687593897.20211008035055.jpg"
$ws.Range("E43").Value = "11100"
$ws.Range("G43").Value = "네이버"
$ws.Range("H43").Value = "일반 - 가격비교 상품"
$ws.Range("I43").Value = ""
$ws.Range("J43").Value = ""
$ws.Range("N43").Value = "장스탠드"

# Row 44 <- source row 43
$ws.Range("B44").Value = "마켓비 장스탠드 엔틱 롱 원룸 장스텐드 조명 등 거실"
$ws.Range("C44").Value = "https://search.shopping.naver.com/gate.nhn?id=82667453108"
$ws.Range("D44").Value = "https://shopping-phinf.pstatic.net/main_8266745/82667453108.jpg"
$ws.Range("E44").Value = "39900"
$ws.Range("G44").Value = "해봄스토어"
$ws.Range("H44").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I44").Value = "마켓비"
$ws.Range("J44").Value = "마켓비"
$ws.Range("N44").Value = "장스탠드"

# Row 45 <- source row 44
$ws.Range("B45").Value = "파파조명 파파 로델리 장스탠드 인테리어스탠드 무드등 조명"
$ws.Range("C45").Value = "https://search.shopping.naver.com/gate.nhn?id=24114394000"
$ws.Range("D45").Value = "https://shopping-phinf.pstatic.net/main_2411439/24114394000.20200910210311.jpg"
$ws.Range("E45").Value = "24146"
$ws.Range("G45").Value = "네이버"
$ws.Range("H45").Value = "일반 - 가격비교 상품"
$ws.Range("I45").Value = "파파조명"
$ws.Range("J45").Value = "베스트조명"
$ws.Range("N45").Value = "장스탠드"

# Row 46 <- source row 45
$ws.Range("B46").Value = "이케아 무드등 LED 장스탠드 조명 식물등 독서등 거실 인테리어"
$ws.Range("C46").Value = "https://search.shopping.naver.com/gate.nhn?id=12543796442"
$ws.Range("D46").Value = "https://shopping-phinf.pstatic.net/main_1254379/12543796442.11.jpg"
$ws.Range("E46").Value = "9000"
$ws.Range("G46").Value = "바이더리빙"
$ws.Range("H46").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I46").Value = "이케아"
$ws.Range("J46").Value = "이케아"
$ws.Range("N46").Value = "장스탠드"

# Row 47 <- source row 46
$ws.Range("B47").Value = "이케아 테르티알 포르소 책상 집게 스탠드 조명 LED 고정 공부 학습용 독서등 작업등"
$ws.Range("C47").Value = "https://search.shopping.naver.com/gate.nhn?id=82252449479"
$ws.Range("D47").Value = "https://shopping-phinf.pstatic.net/main_8225244/82252449479.1.jpg"
$ws.Range("E47").Value = "24290"
$ws.Range("G47").Value = "빅코스트몰"
$ws.Range("H47").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I47").Value = "이케아"
$ws.Range("J47").Value = "이케아"
$ws.Range("N47").Value = "단스탠드"

# Row 48 <- source row 47
$ws.Range("B48").Value = "거실 스탠드조명 북유럽 창의적이고 개성 서재 골드 모던 심플"
$ws.Range("C48").Value = "https://search.shopping.naver.com/gate.nhn?id=17820145800"
$ws.Range("D48").Value = "https://shopping-phinf.pstatic.net/main_1782014/17820145800.20200114013210.jpg"
$ws.Range("E48").Value = "95000"
$ws.Range("G48").Value = "네이버"
$ws.Range("H48").Value = "일반 - 가격비교 상품"
$ws.Range("I48").Value = ""
$ws.Range("J48").Value = ""
$ws.Range("N48").Value = "장스탠드"

# Row 49 <- source row 48
$ws.Range("B49").Value = "[오늘 출발] 이케아 플로어 스탠드 조명 인테리어 레르스타 거실 침실 장스탠드 독서등"
$ws.Range("C49").Value = "https://search.shopping.naver.com/gate.nhn?id=8365733044"
$ws.Range("D49").Value = "https://shopping-phinf.pstatic.net/main_8365733/8365733044.7.jpg"
$ws.Range("E49").Value = "15900"
$ws.Range("G49").Value = "달콩이네 리빙"
$ws.Range("H49").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I49").Value = "이케아"
$ws.Range("J49").Value = "이케아"
$ws.Range("N49").Value = "장스탠드"

# Row 50 <- source row 49
$ws.Range("B50").Value = "활장스탠드 장스탠드 거실등 스탠드조명 플로어스탠드"
$ws.Range("C50").Value = "https://search.shopping.naver.com/gate.nhn?id=29755495858"
$ws.Range("D50").Value = "https://shopping-phinf.pstatic.net/main_2975549/29755495858.20211119232453.jpg"
$ws.Range("E50").Value = "98750"
$ws.Range("G50").Value = "네이버"
$ws.Range("H50").Value = "일반 - 가격비교 상품"
$ws.Range("I50").Value = ""
$ws.Range("J50").Value = ""
$ws.Range("N50").Value = "장스탠드"

# Row 51 <- source row 50
$ws.Range("B51").Value = "마켓비 장스탠드 침실 거실 조명 램프 롱 활장 플로어 독서등 무드등 인테리어"
$ws.Range("C51").Value = "https://search.shopping.naver.com/gate.nhn?id=82685221109"
$ws.Range("D51").Value = "https://shopping-phinf.pstatic.net/main_8268522/82685221109.1.jpg"
$ws.Range("E51").Value = "17400"
$ws.Range("G51").Value = "나누벨"
$ws.Range("H51").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I51").Value = "마켓비"
$ws.Range("J51").Value = "마켓비"
$ws.Range("N51").Value = "장스탠드"

# Row 52 <- source row 51
$ws.Range("B52").Value = "마켓비 단스탠드 라탄 이케아 조명 책상 미니 독서등 침대 침실 우드 엔틱 무드등 갓 램프"
$ws.Range("C52").Value = "https://search.shopping.naver.com/gate.nhn?id=82821838052"
$ws.Range("D52").Value = "https://shopping-phinf.pstatic.net/main_8282183/82821838052.2.jpg"
$ws.Range("E52").Value = "6900"
$ws.Range("G52").Value = "오펜퍼니처"
$ws.Range("H52").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I52").Value = "마켓비"
$ws.Range("J52").Value = "마켓비"
$ws.Range("N52").Value = "단스탠드"

# Row 53 <- source row 52
$ws.Range("B53").Value = "IKEA 이케아 ÅRSTID 오르스티드 플로어스탠드 장스탠드 거실조명 니켈도금, 황동"
$ws.Range("C53").Value = "https://search.shopping.naver.com/gate.nhn?id=81284041579"
$ws.Range("D53").Value = "https://shopping-phinf.pstatic.net/main_8128404/81284041579.jpg"
$ws.Range("E53").Value = "45330"
$ws.Range("G53").Value = "이케아데이"
$ws.Range("H53").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I53").Value = "이케아"
$ws.Range("J53").Value = "이케아"
$ws.Range("N53").Value = "장스탠드"

# Row 54 <- source row 53
$ws.Range("B54").Value = "루이스 폴센 판텔라  플로어 고급형 스탠드 조명 램프 장스탠드  무드등 수입 조명"
$ws.Range("C54").Value = "https://search.shopping.naver.com/gate.nhn?id=29906153623"
$ws.Range("D54").Value = "https://shopping-phinf.pstatic.net/main_2990615/29906153623.jpg"
$ws.Range("E54").Value = "170000"
$ws.Range("G54").Value = "LAPERTA"
$ws.Range("H54").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I54").Value = ""
$ws.Range("J54").Value = ""
$ws.Range("N54").Value = "장스탠드"

# Row 55 <- source row 54
$ws.Range("B55").Value = "앤틱 거실 스탠드 조명 등 바로니스 플로어 장스탠드 LED 전등 결혼선물 집들이선물"
$ws.Range("C55").Value = "https://search.shopping.naver.com/gate.nhn?id=11025477023"
$ws.Range("D55").Value = "https://shopping-phinf.pstatic.net/main_1102547/11025477023.11.jpg"
$ws.Range("E55").Value = "128000"
$ws.Range("G55").Value = "스탠드나라"
$ws.Range("H55").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I55").Value = ""
$ws.Range("J55").Value = "루미앤"
$ws.Range("N55").Value = "장스탠드"

# Row 56 <- source row 55
$ws.Range("B56").Value = "LED독서등 눈보호 LED스탠드 조명 책상 공부 스탠드 미니"
$ws.Range("C56").Value = "https://search.shopping.naver.com/gate.nhn?id=81880952226"
$ws.Range("D56").Value = "https://shopping-phinf.pstatic.net/main_8188095/81880952226.jpg"
$ws.Range("E56").Value = "27800"
$ws.Range("G56").Value = "신바람스토어"
$ws.Range("H56").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I56").Value = "듀플렉스"
$ws.Range("J56").Value = "듀플렉스"
$ws.Range("N56").Value = "단스탠드"

# Row 57 <- source row 56
$ws.Range("B57").Value = "북유럽 거실 침실 스탠드조명 심플 인테리어 거실 소파등 북유럽 심플 스탠드 이케아 장식"
$ws.Range("C57").Value = "https://search.shopping.naver.com/gate.nhn?id=22455317786"
$ws.Range("D57").Value = "https://shopping-phinf.pstatic.net/main_2245531/22455317786.20200411155831.jpg"
$ws.Range("E57").Value = "73210"
$ws.Range("G57").Value = "네이버"
$ws.Range("H57").Value = "일반 - 가격비교 상품"
$ws.Range("I57").Value = ""
$ws.Range("J57").Value = ""
$ws.Range("N57").Value = "장스탠드"

# Row 69 <- source row 70
$ws.Range("B69").Value = "벨라 장스탠드 거실스탠드 플로어 조명 무드등 인테리어스탠드"
$ws.Range("C69").Value = "https://search.shopping.naver.com/gate.nhn?id=28891592179"
$ws.Range("D69").Value = "https://shopping-phinf.pstatic.net/main_2889159/28891592179.20210917221150.jpg"
$ws.Range("E69").Value = "15770"
$ws.Range("G69").Value = "네이버"
$ws.Range("H69").Value = "일반 - 가격비교 상품"
$ws.Range("I69").Value = ""
$ws.Range("J69").Value = ""
$ws.Range("N69").Value = "장스탠드"

# Row 70 <- source row 69
$ws.Range("B70").Value = "이케아 SIMRISHAMN 심리스함 플로어스탠드 전구미포함 조명 장스탠드 거실등"
$ws.Range("C70").Value = "https://search.shopping.naver.com/gate.nhn?id=28793599790"
$ws.Range("D70").Value = "https://shopping-phinf.pstatic.net/main_2879359/28793599790.20211108050557.jpg"
$ws.Range("E70").Value = "98450"
$ws.Range("G70").Value = "네이버"
$ws.Range("H70").Value = "일반 - 가격비교 상품"
$ws.Range("I70").Value = "이케아"
$ws.Range("J70").Value = "이케아"
$ws.Range("N70").Value = "장스탠드"

# Row 88 <- source row 91
$ws.Range("B88").Value = "라루즈 라르고 밝기조절 장스탠드 무드등 식탁 스탠드조명 거실등 인테리어 조명 LED 디밍"
$ws.Range("C88").Value = "https://search.shopping.naver.com/gate.nhn?id=82251549610"
$ws.Range("D88").Value = "https://shopping-phinf.pstatic.net/main_8225154/82251549610.2.jpg"
$ws.Range("E88").Value = "48500"
$ws.Range("G88").Value = "라루즈"
$ws.Range("H88").Value = "일반 - 가격비교 비매칭 일반상품"
$ws.Range("I88").Value = "라루즈"
$ws.Range("J88").Value = ""
$ws.Range("N88").Value = "장스탠드"

# Row 89 <- source row 90
$ws.Range("B89").Value = "파파 튤립 장스탠드 인테리어스탠드 조명"
$ws.Range("C89").Value = "https://search.shopping.naver.com/gate.nhn?id=27109699411"
$ws.Range("D89").Value = "https://shopping-phinf.pstatic.net/main_2710969/27109699411.20210511183001.jpg"
$ws.Range("E89").Value = "26300"
$ws.Range("G89").Value = "네이버"
$ws.Range("H89").Value = "일반 - 가격비교 상품"
$ws.Range("I89").Value = ""
$ws.Range("J89").Value = ""
$ws.Range("N89").Value = "장스탠드"

# Row 90 <- source row 89
$ws.Range("B90").Value = "거실스탠드조명 플로어 램프 ins 북유럽의 창의적인 아이덴티티_ 000015350"
$ws.Range("C90").Value = "https://search.shopping.naver.com/gate.nhn?id=26186427462"
$ws.Range("D90").Value = "https://shopping-phinf.pstatic.net/main_2618642/26186427462.20210228195326.jpg"
$ws.Range("E90").Value = "39690"
$ws.Range("G90").Value = "네이버"
$ws.Range("H90").Value = "일반 - 가격비교 상품"
$ws.Range("I90").Value = ""
$ws.Range("J90").Value = ""
$ws.Range("N90").Value = "장스탠드"

# Row 91 <- source row 88
$ws.Range("B91").Value = "이케아 NYMANE 뉘모네 플로어스탠드3등 조명 거실등 전구미포함"
$ws.Range("C91").Value = "https://search.shopping.naver.com/gate.nhn?id=27063760142"
$ws.Range("D91").Value = "https://shopping-phinf.pstatic.net/main_2706376/27063760142.20210507193639.jpg"
$ws.Range("E91").Value = "61340"
$ws.Range("G91").Value = "네이버"
$ws.Range("H91").Value = "일반 - 가격비교 상품"
$ws.Range("I91").Value = "이케아"
$ws.Range("J91").Value = "이케아"
$ws.Range("N91").Value = "장스탠드"

# Row 93 <- source row 94
$ws.Range("B93").Value = "한샘 한샘몰X스피아노 헤이즈 스탠드 조명SET 플로어 테이블 인테리어스탠드 장스탠드"
$ws.Range("C93").Value = "https://search.shopping.naver.com/gate.nhn?id=26858679216"
$ws.Range("D93").Value = "https://shopping-phinf.pstatic.net/main_2685867/26858679216.20210421193743.jpg"
$ws.Range("E93").Value = "117070"
$ws.Range("G93").Value = "네이버"
$ws.Range("H93").Value = "일반 - 가격비교 상품"
$ws.Range("I93").Value = "한샘"
$ws.Range("J93").Value = ""
$ws.Range("N93").Value = "장스탠드"

# Row 94 <- source row 93
$ws.Range("B94").Value = "북유럽 거실 플로어 깃털스탠드 인테리어 조명"
$ws.Range("C94").Value = "https://search.shopping.naver.com/gate.nhn?id=28876949541"
$ws.Range("D94").Value = "https://shopping-phinf.pstatic.net/main_2887694/28876949541.20210916235401.jpg"
$ws.Range("E94").Value = "39500"
$ws.Range("G94").Value = "네이버"
$ws.Range("H94").Value = "일반 - 가격비교 상품"
$ws.Range("I94").Value = ""
$ws.Range("J94").Value = ""
$ws.Range("N94").Value = "장스탠드"
